$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "nigga" and "dammit" blacklisted word rows (current rows 4 and 5).
# Everything below shifts up automatically.
$ws.Range("A4:D5").EntireRow.Delete() | Out-Null
$ws.Range("B13").Select() | Out-Null
